$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C (pushes old C -> D), preserving per-cell
# styling the same way the original sheet carried it (A:C styled,
# the trailing "translation" column left unstyled on data rows).
$ws.Columns.Item(3).Insert()

# Row 1 (headers): base | en | es | fr
$ws.Cells.Item(1,1).Value = "base"
$ws.Cells.Item(1,2).Value = "en"
$ws.Cells.Item(1,3).Value = "es"
$ws.Cells.Item(1,4).Value = "fr"

# Row 2: fr | plane | avión | avion
$ws.Cells.Item(2,1).Value = "fr"
$ws.Cells.Item(2,2).Value = "plane"
$ws.Cells.Item(2,3).Value = "avión"
$ws.Cells.Item(2,4).Value = "avion"

# Row 3: es | throw | botar | jeter
$ws.Cells.Item(3,1).Value = "es"
$ws.Cells.Item(3,2).Value = "throw"
$ws.Cells.Item(3,3).Value = "botar"
$ws.Cells.Item(3,4).Value = "jeter"

# Row 4: fr | dog | perro | chien
$ws.Cells.Item(4,1).Value = "fr"
$ws.Cells.Item(4,2).Value = "dog"
$ws.Cells.Item(4,3).Value = "perro"
$ws.Cells.Item(4,4).Value = "chien"

# Row 5: en | dog | perro | chien
$ws.Cells.Item(5,1).Value = "en"
$ws.Cells.Item(5,2).Value = "dog"
$ws.Cells.Item(5,3).Value = "perro"
$ws.Cells.Item(5,4).Value = "chien"

# Match the saved selection state (activeCell D5)
$null = $ws.Range("D5").Select()
